$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 167
$ws.Cells.Item(167, 1).Value = 166.0
$ws.Cells.Item(167, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(167, 3).Value = "6:35 PM"
$ws.Cells.Item(167, 4).Value = "FR4529"
$ws.Cells.Item(167, 5).Value = "Oslo"
$ws.Cells.Item(167, 6).Value = "(TRF)"
$ws.Cells.Item(167, 7).Value = "Ryanair "
$ws.Cells.Item(167, 8).Value = "B738"
$ws.Cells.Item(167, 9).Value = "(SP-RSS)"
$ws.Cells.Item(167, 10).Value = "6:20 PM"
$ws.Cells.Item(167, 12).Value = "0 hours, -15 minutes"

# Row 168
$ws.Cells.Item(168, 1).Value = 167.0
$ws.Cells.Item(168, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(168, 3).Value = "7:00 PM"
$ws.Cells.Item(168, 4).Value = "FR1943"
$ws.Cells.Item(168, 5).Value = "Bologna"
$ws.Cells.Item(168, 6).Value = "(BLQ)"
$ws.Cells.Item(168, 7).Value = "Buzz "
$ws.Cells.Item(168, 8).Value = "B38M"
$ws.Cells.Item(168, 9).Value = "(SP-RZE)"
$ws.Cells.Item(168, 10).Value = "6:56 PM"
$ws.Cells.Item(168, 12).Value = "0 hours, -4 minutes"

# Row 169
$ws.Cells.Item(169, 1).Value = 168.0
$ws.Cells.Item(169, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(169, 3).Value = "8:30 PM"
$ws.Cells.Item(169, 4).Value = "FR1939"
$ws.Cells.Item(169, 5).Value = "Gothenburg"
$ws.Cells.Item(169, 6).Value = "(GOT)"
$ws.Cells.Item(169, 7).Value = "Ryanair "
$ws.Cells.Item(169, 8).Value = "B738"
$ws.Cells.Item(169, 9).Value = "(SP-RKD)"
$ws.Cells.Item(169, 10).Value = "8:23 PM"
$ws.Cells.Item(169, 12).Value = "0 hours, -7 minutes"

# Row 170
$ws.Cells.Item(170, 1).Value = 169.0
$ws.Cells.Item(170, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(170, 3).Value = "9:40 PM"
$ws.Cells.Item(170, 4).Value = "FR2019"
$ws.Cells.Item(170, 5).Value = "London"
$ws.Cells.Item(170, 6).Value = "(STN)"
$ws.Cells.Item(170, 7).Value = "Ryanair "
$ws.Cells.Item(170, 8).Value = "B738"
$ws.Cells.Item(170, 9).Value = "(SP-RKP)"
$ws.Cells.Item(170, 10).Value = "9:42 PM"
$ws.Cells.Item(170, 12).Value = "0 hours, 2 minutes"

# Row 171
$ws.Cells.Item(171, 1).Value = 170.0
$ws.Cells.Item(171, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(171, 3).Value = "10:25 PM"
$ws.Cells.Item(171, 4).Value = "FR1055"
$ws.Cells.Item(171, 5).Value = "Brussels"
$ws.Cells.Item(171, 6).Value = "(CRL)"
$ws.Cells.Item(171, 7).Value = "Ryanair "
$ws.Cells.Item(171, 8).Value = "B38M"
$ws.Cells.Item(171, 9).Value = "(SP-RZO)"
$ws.Cells.Item(171, 10).Value = "10:32 PM"
$ws.Cells.Item(171, 12).Value = "0 hours, 7 minutes"

# Row 172
$ws.Cells.Item(172, 1).Value = 171.0
$ws.Cells.Item(172, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(172, 3).Value = "11:00 PM"
$ws.Cells.Item(172, 4).Value = "FR1923"
$ws.Cells.Item(172, 5).Value = "Budapest"
$ws.Cells.Item(172, 6).Value = "(BUD)"
$ws.Cells.Item(172, 7).Value = "Ryanair "
$ws.Cells.Item(172, 8).Value = "B738"
$ws.Cells.Item(172, 9).Value = "(SP-RSV)"
$ws.Cells.Item(172, 10).Value = "10:54 PM"
$ws.Cells.Item(172, 12).Value = "0 hours, -6 minutes"

# Row 173
$ws.Cells.Item(173, 1).Value = 172.0
$ws.Cells.Item(173, 2).Value = "Saturday, Jan 14"
$ws.Cells.Item(173, 3).Value = "11:40 PM"
$ws.Cells.Item(173, 4).Value = "FR6944"
$ws.Cells.Item(173, 5).Value = "Barcelona"
$ws.Cells.Item(173, 6).Value = "(BCN)"
$ws.Cells.Item(173, 7).Value = "Buzz "
$ws.Cells.Item(173, 8).Value = "B38M"
$ws.Cells.Item(173, 9).Value = "(SP-RZG)"
$ws.Cells.Item(173, 10).Value = "11:21 PM"
$ws.Cells.Item(173, 12).Value = "0 hours, -19 minutes"
